# Rebuild the player table (Oyuncu Adı / Pozisyon / Takım) on Sheet1 with the
# new row order and the newly added "Aaron Wiggins" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Amen Thompson",          "PG,SG,SF,PF", "Houston Rockets"),
    @("Anfernee Simons",        "PG,SG",       "Portland Trail Blazers"),
    @("James Harden",           "PG,SG",       "LA Clippers"),
    @("Aaron Wiggins",          "SG,SF",       "Oklahoma City Thunder"),
    @("Anthony Edwards",        "SG,SF",       "Minnesota Timberwolves"),
    @("Jayson Tatum",           "SF,PF",       "Boston Celtics"),
    @("Donte DiVincenzo",       "PG,SG,SF",    "Minnesota Timberwolves"),
    @("RJ Barrett",             "SG,SF,PF",    "Toronto Raptors"),
    @("Ivica Zubac",            "C",           "LA Clippers"),
    @("Giannis Antetokounmpo",  "PF,C",        "Milwaukee Bucks"),
    @("Draymond Green",         "PF,C",        "Golden State Warriors"),
    @("Zion Williamson",        "PF,C",        "New Orleans Pelicans"),
    @("Bradley Beal",           "PG,SG,SF",    "Phoenix Suns"),
    @("Keyonte George",         "PG,SG",       "Utah Jazz"),
    @("Quentin Grimes",         "SG,SF",       "Philadelphia 76ers"),
    @("Jaren Jackson Jr.",      "PF,C",        "Memphis Grizzlies"),
    @("Fred VanVleet",          "PG",          "Houston Rockets"),
    @("Paul George",            "SG,SF,PF",    "Philadelphia 76ers")
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $row = $row + 1
}
